# Weekly driver report update for 2025-04-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths
# (Excel's ColumnWidth property is "characters", which gets re-padded to a
# pixel-snapped "raw" width when written to the sheet XML - offset by
# roughly 0.85-0.92 at this font/DPI. Back the desired raw widths off by a
# safe 0.85 so the saved <col width=.../> lands exactly on the target value.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14 - 0.85
$ws.Columns.Item(5).ColumnWidth = 14 - 0.85
$ws.Columns.Item(6).ColumnWidth = 11 - 0.85
$ws.Columns.Item(7).ColumnWidth = 31 - 0.85
$ws.Columns.Item(8).ColumnWidth = 11 - 0.85
$ws.Columns.Item(9).ColumnWidth = 30 - 0.85
$ws.Columns.Item(10).ColumnWidth = 16 - 0.85

# ---------------------------------------------------------------------------
# "Bad Drivers" table (rows 3-5): reorder the two driver rows and refresh
# the sampled numbers for the new reporting week.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.30.0.6"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 97.59999999999999

$ws.Range("A4").Value = "iwlwifi"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 902
$ws.Range("D4").Value = 97.7

$ws.Range("C5").Value = 903

# ---------------------------------------------------------------------------
# "Good Drivers" table (rows 12-16): the report format changed - it now has
# more columns and the header/data rows no longer carry special styling.
# Wipe the old block (including its formatting) and write the new one.
# ---------------------------------------------------------------------------
$ws.Range("A12:J21").Clear()

$ws.Range("A12").Value = "adapter-driver"
$ws.Range("B12").Value = "good sum"
$ws.Range("C12").Value = "critical sum"
$ws.Range("D12").Value = "warning sum"
$ws.Range("E12").Value = "client count"
$ws.Range("F12").Value = "total sum"
$ws.Range("G12").Value = "adapter"
$ws.Range("H12").Value = "driver"
$ws.Range("I12").Value = "good roaming calculation (%)"
$ws.Range("J12").Value = "driver vintage"

$ws.Range("A13").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.70.3.1"
$ws.Range("B13").Value = 20065
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 29
$ws.Range("F13").Value = 20076
$ws.Range("G13").Value = "intel(r) wi-fi 7 be200 320mhz"
$ws.Range("H13").Value = "23.70.3.1"
$ws.Range("I13").Value = 99.90000000000001
$ws.Range("J13").Formula = "=""2024-08-06"""

$ws.Range("A14").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.60.1.2"
$ws.Range("B14").Value = 47392
$ws.Range("C14").Value = 32
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 96
$ws.Range("F14").Value = 47426
$ws.Range("G14").Value = "intel(r) wi-fi 7 be200 320mhz"
$ws.Range("H14").Value = "23.60.1.2"
$ws.Range("I14").Value = 99.90000000000001
$ws.Range("J14").Formula = "=""2024-06-02"""

$ws.Range("A15").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.40.2.1"
$ws.Range("B15").Value = 27946
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 27946
$ws.Range("G15").Value = "intel(r) wi-fi 7 be200 320mhz"
$ws.Range("H15").Value = "23.40.2.1"
$ws.Range("I15").Value = 100
$ws.Range("J15").Formula = "=""2024-03-30"""

$ws.Range("A16").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.90.0.2"
$ws.Range("B16").Value = 55500
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 89
$ws.Range("F16").Value = 55507
$ws.Range("G16").Value = "intel(r) wi-fi 7 be200 320mhz"
$ws.Range("H16").Value = "23.90.0.2"
$ws.Range("I16").Value = 100
$ws.Range("J16").Formula = "=""2024-09-25"""

# The vintage dates were entered as formulas above purely to dodge Excel's
# automatic date parsing; flatten them back down to literal text values.
$ws.Range("J13:J16").Copy()
$ws.Range("J13:J16").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Trim the now-unused trailing blank rows so the used range matches.
$ws.Rows("17:21").Delete()
